$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) for rows 2-27 changes from 45183 to 45184 (date serial +1 day)
for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
